# Generate Report for Handback
# Two source files were regenerated under new GUID-based names and new
# handoff/handback timestamps. Update every cell (and hyperlink display
# text) that referenced the old file names / timestamps on all three
# sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$oldFile1 = "365764c5-d128-40bc-9cee-edb6cb33f643"
$oldFile2 = "a997d19f-6a67-4018-8d32-d9177a7f1463"
$newFile1 = "0a48022b-3f22-4b64-95fa-057cae1d5fe7"
$newFile2 = "ffff00348e7b-4708-4163-971f-9b27c43aeefa"

$newXlfHash = "1eae068af0547b479fed30e3d824c7431b952718"

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newFile1.md"
$wsOverview.Range("B2").Value = "e2e\$newFile1.md"
$wsOverview.Range("G2").Value = "2016-08-31 21:18:23"

$wsOverview.Range("A3").Value = "$newFile2.md"
$wsOverview.Range("B3").Value = "e2e\$newFile2.md"
$wsOverview.Range("G3").Value = "2016-08-31 21:18:23"

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.TextToDisplay -eq "e2e\$oldFile1.md") {
        $hl.TextToDisplay = "e2e\$newFile1.md"
    } elseif ($hl.TextToDisplay -eq "e2e\$oldFile2.md") {
        $hl.TextToDisplay = "e2e\$newFile2.md"
    }
}

# --- zh-cn sheet -----------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newFile1.md"
$wsZhCn.Range("G2").Value = "$newFile1.$newXlfHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-31 21:18:18"
$wsZhCn.Range("I2").Value = "$newFile1.md"
$wsZhCn.Range("J2").Value = "$newFile1.$newXlfHash.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-31 21:18:36"

$wsZhCn.Range("A3").Value = "$newFile2.md"
$wsZhCn.Range("G3").Value = "$newFile1.$newXlfHash.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-31 21:18:18"
$wsZhCn.Range("I3").Value = "$newFile2.md"
$wsZhCn.Range("J3").Value = "$newFile1.$newXlfHash.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-31 21:18:36"

foreach ($hl in $wsZhCn.Hyperlinks) {
    if ($hl.TextToDisplay -eq "$oldFile1.md") {
        $hl.TextToDisplay = "$newFile1.md"
    } elseif ($hl.TextToDisplay -eq "$oldFile2.md") {
        $hl.TextToDisplay = "$newFile2.md"
    }
}

# --- de-de sheet -----------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newFile1.md"
$wsDeDe.Range("G2").Value = "$newFile1.$newXlfHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-31 21:18:23"
$wsDeDe.Range("I2").Value = "$newFile1.md"
$wsDeDe.Range("J2").Value = "$newFile1.$newXlfHash.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-31 21:18:44"

$wsDeDe.Range("A3").Value = "$newFile2.md"
$wsDeDe.Range("G3").Value = "$newFile1.$newXlfHash.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-31 21:18:23"
$wsDeDe.Range("I3").Value = "$newFile2.md"
$wsDeDe.Range("J3").Value = "$newFile1.$newXlfHash.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-31 21:18:44"

foreach ($hl in $wsDeDe.Hyperlinks) {
    if ($hl.TextToDisplay -eq "$oldFile1.md") {
        $hl.TextToDisplay = "$newFile1.md"
    } elseif ($hl.TextToDisplay -eq "$oldFile2.md") {
        $hl.TextToDisplay = "$newFile2.md"
    }
}
